$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as scraped on 2023-03-04
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.367.38"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.568.06"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.95"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3751"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.08"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3394"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.133"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.93"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.949"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.916"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.567.91"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.87"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06750"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.180"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.50"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.94"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.370.61"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.377"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.701"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.17"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.031"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.45"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.739.55"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.039"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9857"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.995"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.429"
$ws.Range("E36").Value = "  +10.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08458"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02487"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2284"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06468"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.385"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6276"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.15"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.94"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.800"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5902"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.062"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.54"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07326"
$ws.Range("E51").Value = "  +0.82%  "
